# Populate the importer template with its example header/value pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "example"
$ws.Range("A2").Value = "importer"
